# Thai translation edit for "Partner email - document verification failed" template.
# Strategy: use paragraph-scoped Find/Replace (wdReplaceAll=2) so that runs which
# share identical text elsewhere in the document are not accidentally touched,
# and so that replacing a short run (e.g. ", ") does not get merged into an
# adjacent differently-formatted run.

$d = $word.ActiveDocument

function Replace-InRange($range, [string]$find, [string]$replace) {
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# "English" appears twice in the document (hyperlink label in the language
# switcher, and the plain-text heading just below it); both get the same
# translation, so a document-wide replace-all is safe and matches the diff.
Replace-InRange $d.Content "English" "ภาษาอังกฤษ"

$paras = $d.Paragraphs

# Para 1: language switcher line " / Portuguese / French / Thai / Vietnamese / Spanish"
Replace-InRange $paras.Item(1).Range " / Portuguese / French / Thai / Vietnamese / Spanish" " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน"

# Para 5: "Brief" label
Replace-InRange $paras.Item(5).Range "Brief" "บทย่อ"

# Para 6: brief description
Replace-InRange $paras.Item(6).Range "An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io" "อีเมลที่ส่งถึงพันธมิตรในประเทศเป้าหมายที่เอกสารไม่ผ่านการตรวจสอบของเรา โดยมันจะถูกส่งผ่านทาง customer.io"

# Para 8: "Target audience" label
Replace-InRange $paras.Item(8).Range "Target audience" "กลุ่มเป้าหมาย"

# Para 9: target audience description
Replace-InRange $paras.Item(9).Range "Invited partners who submitted wrong/incomplete documents" "พันธมิตรที่ถูกเชิญที่ส่งเอกสารผิดหรือไม่ครบถ้วน"

# Para 12: subject line tail
Replace-InRange $paras.Item(12).Range " — document verification failed " " — การตรวจสอบยืนยันเอกสารล้มเหลว "

# Para 14: hero heading
Replace-InRange $paras.Item(14).Range "Uh oh! Your documents couldn’t be verified" "โอ ไม่นะ! เอกสารของคุณไม่อาจผ่านการตรวจสอบยืนยันได้"

# Para 16: "Hi [PARTNER NAME], "
Replace-InRange $paras.Item(16).Range "Hi " "สวัสดี "
Replace-InRange $paras.Item(16).Range ", " " "

# Para 17: regret message
Replace-InRange $paras.Item(17).Range "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "เราขออภัยที่ต้องแจ้งให้คุณทราบว่า เอกสารของคุณไม่ผ่านกระบวนการตรวจสอบยืนยันของเรา เนื่องจากเราพบปัญหาดังต่อไปนี้: "

# Para 18: document issue bullet
Replace-InRange $paras.Item(18).Range "A copy of your vaccination certificate" "สำเนาใบรับรองการฉีดวัคซีนของคุณ"
Replace-InRange $paras.Item(18).Range ": Document is unclear" ": เอกสารไม่ชัดเจน"

# Para 20: resubmission instructions
Replace-InRange $paras.Item(20).Range "Please resubmit the documents above by " "กรุณายื่นเอกสารข้างต้นอีกครั้งภายในวันที่ "
Replace-InRange $paras.Item(20).Range " so we can proceed with the necessary arrangements." " เพื่อให้เราสามารถดำเนินการตามขั้นตอนที่จำเป็นได้"

# Para 21: live chat / WhatsApp contact line
Replace-InRange $paras.Item(21).Range "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง "
Replace-InRange $paras.Item(21).Range "live chat" "แชทสด"
Replace-InRange $paras.Item(21).Range " or " " หรือทาง "
Replace-InRange $paras.Item(21).Range ". " " "

# Para 22: country manager contact line
Replace-InRange $paras.Item(22).Range "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ "
Replace-InRange $paras.Item(22).Range ", at " " ที่ "
Replace-InRange $paras.Item(22).Range " or " " หรือ "
Replace-InRange $paras.Item(22).Range " (WhatsApp). " " (WhatsApp) "

# Comment text "choose either one" -> Thai. Find does not operate correctly on a
# comment's own Range in this runtime (it silently targets the main document
# story instead), so we must assign the comment range's Text directly.
$comment = $d.Comments.Item(1)
$comment.Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
